$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.625.47'
$ws.Range("E2").Value = '  +2.09%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.888.69'
$ws.Range("E3").Value = '  +0.42%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.88'
$ws.Range("E5").Value = '  +0.89%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9997'
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4915'
$ws.Range("E7").Value = '  -0.18%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2951'
$ws.Range("E8").Value = '  +0.38%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06779'
$ws.Range("E9").Value = '  +2.56%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.887.02'
$ws.Range("E10").Value = '  +0.32%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '17.25'
$ws.Range("E11").Value = '  +3.62%  '
$ws.Range("E12").Value = '  +0.78%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '91.23'
$ws.Range("E13").Value = '  +5.74%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6780'
$ws.Range("E14").Value = '  +1.86%  '
$ws.Range("E15").Value = '  +3.50%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '30.608.80'
$ws.Range("E16").Value = '  +2.17%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000007969'
$ws.Range("E17").Value = '  +2.26%  '
$ws.Range("E18").Value = '  +0.05%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.16'
$ws.Range("E19").Value = '  +3.06%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.131.51'
$ws.Range("E20").Value = '  +0.36%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.002'
$ws.Range("E21").Value = '  +0.35%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.823'
$ws.Range("E22").Value = '  +0.97%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '193.26'
$ws.Range("E23").Value = '  +36.65%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.065'
$ws.Range("E24").Value = '  +3.96%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.335'
$ws.Range("E25").Value = '  +2.81%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '155.23'
$ws.Range("E26").Value = '  +3.32%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.22'
$ws.Range("E27").Value = '  +13.51%  '
$ws.Range("E28").Value = '  +0.38%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.400'
$ws.Range("E29").Value = '  +0.85%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.328'
$ws.Range("E30").Value = '  +3.32%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09077'
$ws.Range("E31").Value = '  +3.72%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.014'
$ws.Range("E32").Value = '  +0.92%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05201'
$ws.Range("E33").Value = '  +3.91%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7563'
$ws.Range("E34").Value = '  +5.34%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.112'
$ws.Range("E35").Value = '  +0.19%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.770'
$ws.Range("E36").Value = '  +3.79%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01834'
$ws.Range("E37").Value = '  +2.34%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.671'
$ws.Range("E38").Value = '  -0.85%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.147'
$ws.Range("E39").Value = '  -0.34%  '
$ws.Range("E40").Value = '  -0.10%  '
$ws.Range("E41").Value = '  +4.64%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '105.04'
$ws.Range("E42").Value = '  +1.70%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.0000'
$ws.Range("E43").Value = '  +0.06%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.741'
$ws.Range("E44").Value = '  +0.42%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.601'
$ws.Range("E45").Value = '  +3.83%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1345'
$ws.Range("E46").Value = '  +6.22%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05854'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.747'
$ws.Range("E48").Value = '  +5.86%  '
$ws.Range("E49").Value = '  +7.15%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.3926'
$ws.Range("E50").Value = '  +4.67%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '33.52'
$ws.Range("E51").Value = '  +2.62%  '
